# Updates crypto price/volume data per upstream refresh (GitHub Actions job)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'25.616.47"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -6.04%  '
$ws.Range('D3').Value = "'1.805.79"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.29%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = "'274.57"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -10.42%  '
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').Value = "'0.5020"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -7.28%  '
$ws.Range('D8').Value = "'0.3495"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -8.38%  '
$ws.Range('E9').Value = '  -4.84%  '
$ws.Range('D10').Value = "'0.06596"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -9.60%  '
$ws.Range('D11').Value = "'19.91"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -10.40%  '
$ws.Range('D12').Value = "'0.8335"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -7.82%  '
$ws.Range('D13').Value = "'0.07764"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.17%  '
$ws.Range('D14').Value = "'1.800.04"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +48.76%  '
$ws.Range('D15').Value = "'5.058"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.56%  '
$ws.Range('D16').Value = "'87.30"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -9.09%  '
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D18').Value = "'13.88"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -6.86%  '
$ws.Range('D19').Value = "'1.001"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').Value = "'0.000007943"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -8.23%  '
$ws.Range('D21').Value = "'25.682.80"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.88%  '
$ws.Range('D22').Value = "'4.715"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -6.74%  '
$ws.Range('D23').Value = "'2.033.25"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +53.64%  '
$ws.Range('D24').Value = "'10.01"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -7.39%  '
$ws.Range('D25').Value = "'6.045"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -7.35%  '
$ws.Range('D26').Value = "'141.87"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.52%  '
$ws.Range('D27').Value = "'2.105"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -8.91%  '
$ws.Range('D28').Value = "'1.652"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.10%  '
$ws.Range('D29').Value = "'16.89"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -8.20%  '
$ws.Range('D30').Value = "'108.09"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.61%  '
$ws.Range('D31').Value = "'4.313"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -11.37%  '
$ws.Range('D32').Value = "'4.181"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -10.63%  '
$ws.Range('D33').Value = "'0.08775"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.81%  '
$ws.Range('D34').Value = "'0.04802"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.52%  '
$ws.Range('D35').Value = "'0.7198"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -13.43%  '
$ws.Range('D36').Value = "'1.124"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -8.25%  '
$ws.Range('D37').Value = "'2.863"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.94%  '
$ws.Range('D38').Value = "'1.000"
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Value = "'3.022"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -9.04%  '
$ws.Range('D40').Value = "'0.01858"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.34%  '
$ws.Range('D41').Value = "'0.5149"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -14.25%  '
$ws.Range('D42').Value = "'2.271"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -15.99%  '
$ws.Range('D43').Value = "'0.9520"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -11.93%  '
$ws.Range('D44').Value = "'113.88"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.19%  '
$ws.Range('D45').Value = "'6.160"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.61%  '
$ws.Range('D46').Value = "'7.982"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -14.11%  '
$ws.Range('D47').Value = "'1.001"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.16%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value = "'0.1376"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -10.24%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').Value = "'0.4533"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -12.12%  '
$ws.Range('D50').Value = "'9.290"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -8.95%  '
$ws.Range('D51').Value = "'35.83"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.41%  '
